$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 328.25
$ws.Cells.Item(29, 10).Value = 987
$ws.Cells.Item(29, 12).Value = 2961
$ws.Cells.Item(29, 14).Value = -3523
$ws.Cells.Item(33, 8).Value = 157.33333
$ws.Cells.Item(33, 9).Value = 220
$ws.Cells.Item(33, 10).Value = 32
$ws.Cells.Item(33, 11).Value = 220
$ws.Cells.Item(33, 12).Value = 32
$ws.Cells.Item(33, 13).Value = 9
$ws.Cells.Item(33, 14).Value = -490
$ws.Cells.Item(40, 8).Value = 3000
$ws.Cells.Item(40, 9).Value = 2440
$ws.Cells.Item(40, 10).Value = 3933.3333
$ws.Cells.Item(40, 11).Value = 2440
$ws.Cells.Item(40, 12).Value = 3933.3333
$ws.Cells.Item(40, 13).Value = -2265
$ws.Cells.Item(40, 14).Value = -4283.3333
$ws.Cells.Item(42, 8).Value = 285.75
$ws.Cells.Item(42, 9).Value = 131
$ws.Cells.Item(42, 10).Value = 750
$ws.Cells.Item(42, 11).Value = 393
$ws.Cells.Item(42, 12).Value = 2250
$ws.Cells.Item(42, 13).Value = -163
$ws.Cells.Item(42, 14).Value = -2710
$ws.Cells.Item(51, 8).Value = 2696
$ws.Cells.Item(51, 9).Value = 2250
$ws.Cells.Item(51, 10).Value = 2993.3333
$ws.Cells.Item(51, 11).Value = 2250
$ws.Cells.Item(51, 12).Value = 2993.3333
$ws.Cells.Item(51, 13).Value = -1766
$ws.Cells.Item(51, 14).Value = -3961.3333
$ws.Cells.Item(62, 8).Value = 2663.3823
$ws.Cells.Item(62, 9).Value = 2362.4
$ws.Cells.Item(62, 10).Value = 3499.4443
$ws.Cells.Item(62, 11).Value = 2362.4
$ws.Cells.Item(62, 12).Value = 3499.4443
$ws.Cells.Item(62, 13).Value = -1738.4
$ws.Cells.Item(62, 14).Value = -4747.4443
$ws.Cells.Item(65, 8).Value = 2663.3823
$ws.Cells.Item(65, 9).Value = 2362.4
$ws.Cells.Item(65, 10).Value = 3499.4443
$ws.Cells.Item(65, 11).Value = 11812
$ws.Cells.Item(65, 12).Value = 17497.2215
$ws.Cells.Item(65, 13).Value = -8692
$ws.Cells.Item(65, 14).Value = -23737.2215

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1470.8975
$ws.Cells.Item(132, 9).Value = 1118.2572
$ws.Cells.Item(132, 10).Value = 4556.5
$ws.Cells.Item(132, 11).Value = 3354.7716
$ws.Cells.Item(132, 12).Value = 13669.5
$ws.Cells.Item(132, 13).Value = -824.7716
$ws.Cells.Item(132, 14).Value = -18729.5
$ws.Cells.Item(138, 8).Value = 36980
$ws.Cells.Item(138, 10).Value = 36980
$ws.Cells.Item(138, 12).Value = 36980
$ws.Cells.Item(138, 14).Value = -47260

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 437.125
$ws.Cells.Item(67, 8).Value = 437.125
$ws.Cells.Item(115, 8).Value = 25155.25
$ws.Cells.Item(115, 10).Value = 16666.666
$ws.Cells.Item(115, 12).Value = 16666.666
$ws.Cells.Item(115, 14).Value = -19800.666

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2439.4443
$ws.Cells.Item(132, 9).Value = 1710.8572
$ws.Cells.Item(132, 11).Value = 5132.571599999999
$ws.Cells.Item(132, 13).Value = -2602.571599999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 409.7
$ws.Cells.Item(68, 9).Value = 350.5
$ws.Cells.Item(68, 10).Value = 449.16666
$ws.Cells.Item(68, 11).Value = 1051.5
$ws.Cells.Item(68, 12).Value = 1347.49998
$ws.Cells.Item(68, 13).Value = -240.5
$ws.Cells.Item(68, 14).Value = -2969.49998
$ws.Cells.Item(71, 8).Value = 409.7
$ws.Cells.Item(71, 9).Value = 350.5
$ws.Cells.Item(71, 10).Value = 449.16666
$ws.Cells.Item(71, 11).Value = 3154.5
$ws.Cells.Item(71, 12).Value = 4042.49994
$ws.Cells.Item(71, 13).Value = 901.5
$ws.Cells.Item(71, 14).Value = -12154.49994
$ws.Cells.Item(87, 8).Value = 35087.7
$ws.Cells.Item(87, 9).Value = 8000
$ws.Cells.Item(87, 10).Value = 35430.582
$ws.Cells.Item(87, 11).Value = 24000
$ws.Cells.Item(87, 12).Value = 106291.746
$ws.Cells.Item(87, 13).Value = -22752
$ws.Cells.Item(87, 14).Value = -108787.746
$ws.Cells.Item(90, 8).Value = 35087.7
$ws.Cells.Item(90, 9).Value = 8000
$ws.Cells.Item(90, 10).Value = 35430.582
$ws.Cells.Item(90, 11).Value = 72000
$ws.Cells.Item(90, 12).Value = 318875.238
$ws.Cells.Item(90, 13).Value = -65760
$ws.Cells.Item(90, 14).Value = -331355.238
$ws.Cells.Item(93, 8).Value = 4384.136
$ws.Cells.Item(93, 9).Value = 1881
$ws.Cells.Item(93, 10).Value = 4940.3887
$ws.Cells.Item(93, 11).Value = 5643
$ws.Cells.Item(93, 12).Value = 14821.1661
$ws.Cells.Item(93, 13).Value = -3771
$ws.Cells.Item(93, 14).Value = -18565.1661
$ws.Cells.Item(96, 8).Value = 35955.555
$ws.Cells.Item(96, 10).Value = 35955.555
$ws.Cells.Item(96, 12).Value = 107866.665
$ws.Cells.Item(96, 14).Value = -111984.665
$ws.Cells.Item(100, 8).Value = 2704
$ws.Cells.Item(100, 10).Value = 2704
$ws.Cells.Item(100, 12).Value = 8112
$ws.Cells.Item(100, 14).Value = -9734
$ws.Cells.Item(106, 8).Value = 4400
$ws.Cells.Item(106, 9).Value = 2000
$ws.Cells.Item(106, 11).Value = 6000
$ws.Cells.Item(106, 13).Value = -5054
$ws.Cells.Item(109, 8).Value = 3157.3057
$ws.Cells.Item(109, 9).Value = 304.4
$ws.Cells.Item(109, 10).Value = 4254.577
$ws.Cells.Item(109, 11).Value = 913.1999999999999
$ws.Cells.Item(109, 12).Value = 12763.731
$ws.Cells.Item(109, 13).Value = 126.8000000000001
$ws.Cells.Item(109, 14).Value = -14843.731
$ws.Cells.Item(112, 8).Value = 3411.25
$ws.Cells.Item(112, 9).Value = 1590
$ws.Cells.Item(112, 10).Value = 3671.4285
$ws.Cells.Item(112, 11).Value = 4770
$ws.Cells.Item(112, 12).Value = 11014.2855
$ws.Cells.Item(112, 13).Value = -3662
$ws.Cells.Item(112, 14).Value = -13230.2855

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7275
$ws.Cells.Item(132, 9).Value = 9073.5
$ws.Cells.Item(132, 10).Value = 3678
$ws.Cells.Item(132, 11).Value = 27220.5
$ws.Cells.Item(132, 12).Value = 11034
$ws.Cells.Item(132, 13).Value = -24690.5
$ws.Cells.Item(132, 14).Value = -16094
$ws.Cells.Item(138, 8).Value = 20195.666
$ws.Cells.Item(138, 10).Value = 20195.666
$ws.Cells.Item(138, 12).Value = 20195.666
$ws.Cells.Item(138, 14).Value = -30475.666

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4200.75
$ws.Cells.Item(40, 9).Value = 3123.1538
$ws.Cells.Item(40, 10).Value = 5134.6665
$ws.Cells.Item(40, 11).Value = 3123.1538
$ws.Cells.Item(40, 12).Value = 5134.6665
$ws.Cells.Item(40, 13).Value = -2987.1538
$ws.Cells.Item(40, 14).Value = -5406.6665
$ws.Cells.Item(122, 8).Value = 3697.3333
$ws.Cells.Item(122, 9).Value = 3551
$ws.Cells.Item(122, 11).Value = 10653
$ws.Cells.Item(122, 13).Value = -8203
$ws.Cells.Item(132, 8).Value = 1206.2297
$ws.Cells.Item(132, 9).Value = 1060.9474
$ws.Cells.Item(132, 11).Value = 3182.8422
$ws.Cells.Item(132, 13).Value = -652.8422
$ws.Cells.Item(135, 8).Value = 27450
$ws.Cells.Item(135, 10).Value = 27450
$ws.Cells.Item(135, 12).Value = 27450
$ws.Cells.Item(135, 14).Value = -37590

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2423.238
$ws.Cells.Item(126, 9).Value = 2449.875
$ws.Cells.Item(126, 10).Value = 2406.8462
$ws.Cells.Item(126, 11).Value = 7349.625
$ws.Cells.Item(126, 12).Value = 7220.5386
$ws.Cells.Item(126, 13).Value = -4879.625
$ws.Cells.Item(126, 14).Value = -12160.5386
